# Add Coeff of Var input data
#
# Inserts a new "CV-CDF" section (3 data rows + 1 blank separator row) at the
# top of the Model_Input table, right after the header row, pushing the
# existing sections (ESTRESS, Groundwater, landocean, OUTPUTQ_LTM, population,
# RiverBasin) down by four rows. Also widens column E so the new, longer
# description text fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model_Input")

# Make room: insert 4 whole rows starting at row 3 (3 new data rows + 1 blank
# separator row, matching the spacing pattern already used between sections).
$ws.Range("A3:A6").EntireRow.Insert()

# Row 3: CV-CDF section header
$ws.Cells.Item(3, 1).Value = "CV-CDF"
$ws.Cells.Item(3, 1).Font.Bold = $true
$ws.Cells.Item(3, 1).Font.Size = 12
$ws.Cells.Item(3, 2).Value = "WBM_TerraClimate2000-2020_Q_DIST_CV_InterAnnual_6min"
$ws.Cells.Item(3, 3).Value = "GeoTiff"
$ws.Cells.Item(3, 4).Value = "unitless"
$ws.Cells.Item(3, 5).Value = "Inter-annual Coefficient of Variation for annual discharges across years 2000-2020"

# Row 4
$ws.Cells.Item(4, 2).Value = "WBM_TerraClimate2000-2020_Q_DIST_CV_IntraAnnual_6min"
$ws.Cells.Item(4, 3).Value = "GeoTiff"
$ws.Cells.Item(4, 4).Value = "unitless"
$ws.Cells.Item(4, 5).Value = "Intra-annual Coefficient of Variation across long-term mean monthly discharges for years 2000-2020"

# Row 5
$ws.Cells.Item(5, 2).Value = "WBM_TerraClimate2000-2020_Q_DIST_meanCV_IntraAnnual_6min"
$ws.Cells.Item(5, 3).Value = "GeoTiff"
$ws.Cells.Item(5, 4).Value = "unitless"
$ws.Cells.Item(5, 5).Value = "Intra-annual Coefficient of Variation of discharges across all months and years from 2000-2020"

# Row 6 stays blank - separator row, consistent with the rest of the sheet.

# Widen column E to fit the new, longer description text (closest width this
# engine's quantized ColumnWidth setter can reach to the target 87.6640625).
$ws.Columns.Item(5).ColumnWidth = 86.8
